$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTop = 44260.52179451851
$newMid = 44260.50046984954
$newBot = 44260.47913657407

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newTop
}
for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $newMid
}
for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $newBot
}
